$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced to Text
# format first, otherwise Excel auto-converts them (e.g. "1.00" -> 1),
# which would not match the source data (these columns store text).
$textCells = @("D5", "D6", "D9", "D12", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D29", "D30", "D31", "D32", "D33", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "64.480.04"
$ws.Range("E2").Value = "  +2.00%  "
$ws.Range("D3").Value = "3.458.46"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "573.63"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "159.71"
$ws.Range("E6").Value = "  +3.63%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "3.467.00"
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("D9").Value = "0.584"
$ws.Range("E9").Value = "  +10.83%  "
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("E11").Value = "  +4.32%  "
$ws.Range("D12").Value = "0.446"
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("D13").Value = "4.052.87"
$ws.Range("E13").Value = "  +2.40%  "
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("E15").Value = "  +6.12%  "
$ws.Range("D16").Value = "28.73"
$ws.Range("E16").Value = "  +5.94%  "
$ws.Range("D17").Value = "64.522.20"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "3.467.26"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("D19").Value = "6.41"
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("D20").Value = "14.38"
$ws.Range("E20").Value = "  +3.14%  "
$ws.Range("D21").Value = "386.75"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").Value = "8.17"
$ws.Range("E22").Value = "  -3.16%  "
$ws.Range("D23").Value = "73.24"
$ws.Range("E23").Value = "  +3.78%  "
$ws.Range("D24").Value = "0.544"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  +17.74%  "
$ws.Range("D27").Value = "9.52"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "6.18"
$ws.Range("E30").Value = "  +9.92%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "1.44"
$ws.Range("E31").Value = "  +10.38%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").Value = "6.66"
$ws.Range("E32").Value = "  +3.24%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "2.03"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("D36").Value = "7.04"
$ws.Range("E36").Value = "  +4.20%  "
$ws.Range("D37").Value = "1.50"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").Value = "160.73"
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.0775"
$ws.Range("E39").Value = "  +3.73%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.88"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").Value = "27.25"
$ws.Range("E41").Value = "  -1.32%  "
$ws.Range("D42").Value = "2.910.87"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").Value = "0.0319"
$ws.Range("E43").Value = "  -3.11%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "4.46"
$ws.Range("E44").Value = "  +3.88%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "42.17"
$ws.Range("E45").Value = "  +2.87%  "
$ws.Range("D46").Value = "0.772"
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("D47").Value = "23.80"
$ws.Range("E47").Value = "  +7.83%  "
$ws.Range("D48").Value = "1.09"
$ws.Range("E48").Value = "  +4.07%  "
$ws.Range("E49").Value = "  +17.52%  "
$ws.Range("E50").Value = "  +4.42%  "
$ws.Range("D51").Value = "0.863"
$ws.Range("E51").Value = "  +6.39%  "
